$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Mã phòng TNTH" (room code) column D contained a stale UUID for both
# equipment rows. Update it to the corrected UUID for both Macbook Pro 1
# and Macbook Pro 2 entries.
$ws.Range("D2").Value = "6b38ddb6-1cdc-5f7b-6efa-d9c911cf4972"
$ws.Range("D3").Value = "6b38ddb6-1cdc-5f7b-6efa-d9c911cf4972"

# Leave the selection on the last-edited cell, as Excel would after typing.
$ws.Range("D3").Select()
